$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "dev" in G1, matching style of the other header cells (F1)
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "dev"

# Populate deviation column G (y - mu = C - F) for each data row
$ws.Range("G2").Value = -5.603313142733896
$ws.Range("G3").Value = 13.44178131332311
$ws.Range("G4").Value = -0.4109123194755568
$ws.Range("G5").Value = -5.847117855493366
$ws.Range("G6").Value = -3.445840932899443
$ws.Range("G7").Value = -4.373669846551962
$ws.Range("G8").Value = -5.789761607707192
$ws.Range("G9").Value = 18.23958771180645
$ws.Range("G10").Value = 15.18559232557055
$ws.Range("G11").Value = 5.145006423519618
$ws.Range("G12").Value = -1.529305952095314
$ws.Range("G13").Value = -10.90580348223946
$ws.Range("G14").Value = 5.692681614088414
$ws.Range("G15").Value = -6.725968792313893
$ws.Range("G16").Value = -1.741455315965538
$ws.Range("G17").Value = 1.749274073039686
$ws.Range("G18").Value = -3.701725396732314
$ws.Range("G19").Value = 0.3376440857239515
$ws.Range("G20").Value = -2.300460305191415
$ws.Range("G21").Value = 17.32180131382403
$ws.Range("G22").Value = 17.56528025292826
$ws.Range("G23").Value = 17.30731329899686
$ws.Range("G24").Value = 5.523193298495954
$ws.Range("G25").Value = -1.493256697556369
$ws.Range("G26").Value = -2.348742132203867
$ws.Range("G27").Value = -2.357624903495577
$ws.Range("G28").Value = 4.589348655310374
$ws.Range("G29").Value = -7.611607599685385
$ws.Range("G30").Value = -1.096309997399757
$ws.Range("G31").Value = -12.54295471568719
$ws.Range("G32").Value = 5.806808429897416
$ws.Range("G33").Value = -2.229316744888308
$ws.Range("G34").Value = -13.79736693294724
$ws.Range("G35").Value = -7.247221382652
$ws.Range("G36").Value = -0.1143572385845459
$ws.Range("G37").Value = -6.271950829362254
$ws.Range("G38").Value = -1.583546734511032
$ws.Range("G39").Value = 6.069930888854458
$ws.Range("G40").Value = -0.3655451689015479
$ws.Range("G41").Value = -20.29510415394371
$ws.Range("G42").Value = 9.86518049060146
$ws.Range("G43").Value = -2.231821832189908
$ws.Range("G44").Value = 3.278374871735487
$ws.Range("G45").Value = -9.537505441904216
$ws.Range("G46").Value = -6.112103105200845
$ws.Range("G47").Value = 7.117969553455936
$ws.Range("G48").Value = -15.57292348049555
$ws.Range("G49").Value = 4.23158867196436
$ws.Range("G50").Value = 13.48211321252276
$ws.Range("G51").Value = -6.76587644465107
